$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore DB feature implemented: mark "I want to export and reimport DB" as DONE
$ws.Range("C12").Value = "DONE"
# Next sprint-4 story moves from NOT STARTED to IN PROGRESS
$ws.Range("C13").Value = "IN PROGRESS"

# Apply AutoFilter on the Status column (col 3) to show only IN PROGRESS and NOT STARTED
$ws.Range("A1:C15").AutoFilter(3, @("IN PROGRESS", "NOT STARTED"), 7)
